$wb = $excel.ActiveWorkbook

$wsObjType = $wb.Worksheets.Item("ObjType")
$wsObjType.Activate()
$wsObjType.Range("A22").Select()

$ws = $wb.Worksheets.Item("Common")
$ws.Activate()

$ws.Range("A58").Value = "EWorldType"
$ws.Range("B58").Value = "NORMAL"
$ws.Range("C58").Value = 0

$ws.Range("A59").Value = "EWorldType"
$ws.Range("B59").Value = "DARK"
$ws.Range("C59").Value = 1

$ws.Range("A60").Value = "EWorldType"
$ws.Range("B60").Value = "MASTER"
$ws.Range("C60").Value = 2

$ws.Range("A61").Value = "EWorldStageType"
$ws.Range("B61").Value = "STAGE"
$ws.Range("C61").Value = 0

$ws.Range("A62").Value = "EWorldStageType"
$ws.Range("B62").Value = "REWARD"
$ws.Range("C62").Value = 1

$ws.Range("A63").Value = "EWorldStageType"
$ws.Range("B63").Value = "VILLAGE"
$ws.Range("C63").Value = 2

$ws.Range("D60:E64").Select()
